# NSMB workbook update
# - fix map movements after 5-1 (lost 17 frames to luck manipulation)
# - record new "fail"/"Fail"/"WIN" markers and move/frame-count columns (J/K/I/L)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Corrected B/H-column timings (luck manipulation fix) ---
$ws.Range("B39").Value = 15307
$ws.Range("B40").Value = 15573
$ws.Range("B42").Value = 15722
$ws.Range("H42").Value = 15781

# --- New outcome markers (K column) + move counters (J column) ---
# Written in an order that matches how the shared-string table grows:
# "fail", "Fail", "WIN", "15290 possible".
$ws.Range("J35").Value = 76
$ws.Range("K35").Value = "fail"

$ws.Range("J36").Value = 77
$ws.Range("K36").Value = "fail"

$ws.Range("J37").Value = 78

$ws.Range("J38").Value = 79
$ws.Range("K38").Value = "fail"

$ws.Range("J39").Value = 80

$ws.Range("J40").Value = 81

$ws.Range("J41").Value = 82

$ws.Range("J42").Value = 83
$ws.Range("K42").Value = "Fail"

$ws.Range("J43").Value = 84

$ws.Range("J44").Value = 85
$ws.Range("K44").Value = "Fail"

$ws.Range("J45").Value = 86

$ws.Range("J46").Value = 87

$ws.Range("J47").Value = 88

$ws.Range("J48").Value = 89

$ws.Range("J49").Value = 90

$ws.Range("J50").Value = 91

$ws.Range("J51").Value = 92
$ws.Range("K51").Value = "fail"

$ws.Range("I52").Value = 95090
$ws.Range("J52").Value = 93
$ws.Range("K52").Value = "WIN"
$ws.Range("L52").Value = 94790

# "15290 possible" is introduced last in the shared-string table.
$ws.Range("I39").Value = "15290 possible"

# --- View state: move the active cell/selection to K25 ---
$ws.Range("K25").Select()

$wb.Save()
